$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "28.952.02"
$ws.Range("E2").Value = "  +1.03%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.887.61"
$ws.Range("E3").Value = "  +0.64%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.002"
$ws.Range("E4").Value = "  -0.32%  "
Set-TextValue $ws.Cells.Item(5, 4) "325.55"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  -0.35%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.4585"
$ws.Range("E7").Value = "  +1.04%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3899"
$ws.Range("E8").Value = "  +1.49%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.07858"
$ws.Range("E9").Value = "  +0.46%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.9879"
$ws.Range("E10").Value = "  -0.26%  "
Set-TextValue $ws.Cells.Item(11, 4) "21.91"
$ws.Range("E11").Value = "  +1.89%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.908.86"
$ws.Range("E12").Value = "  +0.66%  "
Set-TextValue $ws.Cells.Item(13, 4) "7.024"
$ws.Range("E13").Value = "  +1.52%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.693"
$ws.Range("E14").Value = "  +0.91%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.06951"
$ws.Range("E15").Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(16, 4) "88.17"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  -0.36%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.000009990"
$ws.Range("E18").Value = "  +0.31%  "
Set-TextValue $ws.Cells.Item(19, 4) "17.03"
$ws.Range("E19").Value = "  +1.83%  "
Set-TextValue $ws.Cells.Item(20, 4) "1.002"
$ws.Range("E20").Value = "  -0.19%  "
Set-TextValue $ws.Cells.Item(21, 4) "28.975.89"
$ws.Range("E21").Value = "  +1.13%  "
Set-TextValue $ws.Cells.Item(22, 4) "5.286"
$ws.Range("E22").Value = "  +0.52%  "
Set-TextValue $ws.Cells.Item(23, 4) "10.97"
$ws.Range("E23").Value = "  +0.85%  "
Set-TextValue $ws.Cells.Item(24, 4) "2.114.92"
$ws.Range("E24").Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.053"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  +0.83%  "
Set-TextValue $ws.Cells.Item(27, 4) "19.32"
$ws.Range("E27").Value = "  +0.94%  "
Set-TextValue $ws.Cells.Item(28, 4) "5.987"
$ws.Range("E28").Value = "  +5.77%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.931"
$ws.Range("E29").Value = "  +2.73%  "
Set-TextValue $ws.Cells.Item(30, 4) "117.64"
$ws.Range("E30").Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.09341"
$ws.Range("E31").Value = "  +0.74%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.9079"
$ws.Range("E32").Value = "  +0.54%  "
Set-TextValue $ws.Cells.Item(33, 4) "5.286"
$ws.Range("E33").Value = "  +0.25%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.330"
$ws.Range("E34").Value = "  +0.75%  "
Set-TextValue $ws.Cells.Item(35, 4) "3.259"
$ws.Range("E35").Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(36, 4) "1.207"
$ws.Range("E36").Value = "  +4.86%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.05777"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("E39").Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(40, 4) "7.661"
$ws.Range("E40").Value = "  +0.67%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.5685"
$ws.Range("E41").Value = "  +2.40%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.1772"
$ws.Range("E42").Value = "  +0.39%  "
Set-TextValue $ws.Cells.Item(43, 4) "9.756"
$ws.Range("E43").Value = "  +1.81%  "
Set-TextValue $ws.Cells.Item(44, 4) "2.270"
$ws.Range("E44").Value = "  +5.43%  "
Set-TextValue $ws.Cells.Item(45, 4) "11.91"
$ws.Range("E45").Value = "  +3.53%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.5369"
$ws.Range("E46").Value = "  +2.38%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.07034"
$ws.Range("E47").Value = "  -1.75%  "
Set-TextValue $ws.Cells.Item(49, 4) "112.93"
$ws.Range("E49").Value = "  +1.40%  "
Set-TextValue $ws.Cells.Item(50, 4) "2.534"
$ws.Range("E50").Value = "  +3.99%  "
$ws.Range("E51").Value = "  -3.31%  "
